$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tail of the sheet (A78:A82) currently holds the standalone "冊" (volume) name
# list ("第二冊".."第六冊"). Move it down 5 rows to A83:A87 to make room for five
# new data rows, mirroring the author's edit.
$ws.Range("A78:A82").Cut($ws.Range("A83:A87"))

# --- Row 78 ---------------------------------------------------------------
$ws.Range("C78").Value = "第2單元"
$ws.Range("A78").Value = "總複習"
$ws.Range("B78").Value = "第一冊"
$ws.Range("D78").Value = "P27"
$ws.Range("E78").Value = "https://view.genially.com/65ce270ed2092a0013c9dab9/guide-2yuan"

# --- Row 79 ---------------------------------------------------------------
$ws.Range("C79").Value = "第3單元"
$ws.Range("A79").Value = "總複習"
$ws.Range("B79").Value = "第一冊"
$ws.Range("D79").Value = "P33"
$ws.Range("E79").Value = "https://view.genially.com/669cd1dd079a5f293384be18/learning-experience-didactic-unit-3yuanyuan"

# --- Row 80 ---------------------------------------------------------------
$ws.Range("C80").Value = "第4單元"
$ws.Range("A80").Value = "總複習"
$ws.Range("B80").Value = "第一冊"
$ws.Range("D80").Value = "P46"
$ws.Range("E80").Value = "https://view.genially.com/669e8c4097ad7b0b01e033ab/interactive-content-4yuan-yuan"

# --- Row 81 ---------------------------------------------------------------
$ws.Range("C81").Value = "1-2"
$ws.Range("A81").Value = "八體"
$ws.Range("B81").Value = "B1"
$ws.Range("D81").Value = "P35"
$ws.Range("E81").Value = "https://view.genially.com/669fb5d47dadb63f096d4905/guide-b1-1-2"

# --- Row 82 ---------------------------------------------------------------
$ws.Range("C82").Value = "1-3"
$ws.Range("A82").Value = "八體"
$ws.Range("B82").Value = "B1"
$ws.Range("D82").Value = "P40"
$ws.Range("E82").Value = "https://view.genially.com/669fbec44b779ff586d1254f/presentation-b1-1-3"

# The unit column (C) on rows 78-80 keeps the plain/general format (no explicit
# style override) -- copy the unstyled format from A1 onto them.
$ws.Range("A1").Copy()
$ws.Range("C78:C80").PasteSpecial(-4122)

# Column E holds hyperlinks styled like the rest of the sheet's link cells --
# copy that cell format (style index 2) onto the new link cells.
$ws.Range("E77").Copy()
$ws.Range("E78:E82").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Wire up the actual hyperlinks (adds entries to sheet1.xml.rels too).
$ws.Hyperlinks.Add($ws.Range("E78"), "https://view.genially.com/65ce270ed2092a0013c9dab9/guide-2yuan") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E79"), "https://view.genially.com/669cd1dd079a5f293384be18/learning-experience-didactic-unit-3yuanyuan") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E80"), "https://view.genially.com/669e8c4097ad7b0b01e033ab/interactive-content-4yuan-yuan") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E81"), "https://view.genially.com/669fb5d47dadb63f096d4905/guide-b1-1-2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E82"), "https://view.genially.com/669fbec44b779ff586d1254f/presentation-b1-1-3") | Out-Null

# Re-apply the link-style format after Hyperlinks.Add (it can introduce its own
# style entry) so the cells keep reusing the sheet's existing hyperlink style.
$ws.Range("E77").Copy()
$ws.Range("E78:E82").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll / selection bookkeeping to match the final view state.
$ws.Range("E83").Select()
